# Update gh-pages to output generated at 2291077
$wb = $excel.ActiveWorkbook

# Sheets: 展览 (Exhibition), 演出 (Performance), 本地生活 (Local life), 全部类型 (All types)
$sheetNames = @("展览", "演出", "本地生活", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("J1").Value = "Cover"
}

# Sheets with data rows (展览, 全部类型) get updated counts + cover image links
$dataSheetNames = @("展览", "全部类型")

foreach ($name in $dataSheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 227
    $ws.Range("J2").Value = "//i1.hdslb.com/bfs/openplatform/202311/lP5IkqWn1699431829470.jpeg"

    $ws.Range("F3").Value = 236
    $ws.Range("J3").Value = "//i1.hdslb.com/bfs/openplatform/202312/ee5hLUN61702276208812.jpeg"
}
